$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B1").Value = 1751037614.045149
$ws.Range("B2").Value = 311878035
$ws.Range("C2").Value = 311877623
$ws.Range("B3").Value = 311878035
$ws.Range("C3").Value = 311877623
$ws.Range("B4").Value = 311876387
$ws.Range("B5").Value = 311870207
$ws.Range("B7").Value = 311867735
$ws.Range("B10").Value = 311867323
$ws.Range("B11").Value = 408296512.08
$ws.Range("B12").Value = 407128632.3996
$ws.Range("B13").Value = 311867735
$ws.Range("B14").Value = 311867735
$ws.Range("B17").Value = 311867735
$ws.Range("B18").Value = 407123960.08
$ws.Range("B20").Value = 407129032.84
$ws.Range("B21").Value = 407128637.434
$ws.Range("B23").Value = 369884507.0368
$ws.Range("B25").Value = 502390362.0508
$ws.Range("B26").Value = 369884858.4719999
$ws.Range("B27").Value = 502389936.1276
$ws.Range("B28").Value = 311867323
$ws.Range("B29").Value = 311867323
$ws.Range("B30").Value = 311867323
$ws.Range("B33").Value = 407129032.84
$ws.Range("B34").Value = 311867323
$ws.Range("B35").Value = 311867323
$ws.Range("B37").Value = 311867323
$ws.Range("B38").Value = 407129048.16
$ws.Range("B40").Value = 408304064.04
$ws.Range("B41").Value = 407128627.04
$ws.Range("B42").Value = 311867323
$ws.Range("B43").Value = 311867735
$ws.Range("B44").Value = 311868147
$ws.Range("B46").Value = 311867323
$ws.Range("B50").Value = 407129839.9741279
$ws.Range("B51").Value = 311867323
$ws.Range("B52").Value = 311867323
$ws.Range("B53").Value = 311867323
$ws.Range("B54").Value = 407129037.304
$ws.Range("B56").Value = 407129039.44
$ws.Range("B58").Value = 311867323
$ws.Range("B60").Value = 369954471.84
$ws.Range("B61").Value = 407127197.44
$ws.Range("B62").Value = 407128634.3092
$ws.Range("B63").Value = 311867323
$ws.Range("B66").Value = 311867735
$ws.Range("B67").Value = 311867323
$ws.Range("B69").Value = 311867735
$ws.Range("B70").Value = 311867735
$ws.Range("B73").Value = 502390263.5903448
$ws.Range("B74").Value = 311868147
$ws.Range("B75").Value = 407128623.84
$ws.Range("B76").Value = 311867323
$ws.Range("B78").Value = 311867323
$ws.Range("B79").Value = 407129032.84
$ws.Range("B81").Value = 311867323
$ws.Range("B82").Value = 311867323
$ws.Range("B84").Value = 407899472.84
$ws.Range("B85").Value = 311868147
$ws.Range("B86").Value = 407129205.52
$ws.Range("B87").Value = 407128628.6
$ws.Range("B88").Value = 465145730.44
$ws.Range("B91").Value = 311868147
$ws.Range("B93").Value = 311867323
$ws.Range("B94").Value = 407123960.08
$ws.Range("B95").Value = 311867323
$ws.Range("B96").Value = 311868147
